$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = "-"
$ws.Range("E4").Value = "-"

# Row 8
$ws.Range("C8").Value = "-"
$ws.Range("E8").Value = "-"

# Row 11
$ws.Range("C11").Value = "-"

# Row 12
$ws.Range("C12").Value = "-"

# Row 14
$ws.Range("C14").Value = "-"

# Row 15
$ws.Range("C15").Value = "-"

# Row 18
$ws.Range("B18").Value = "[Emerson-Eletrônica Básica-2NA, Allan Cupertino-Instalções Elétricas-2NA]"
$ws.Range("C18").Value = "[Cláudio-Tecnologia da Soldagem-2NA, Paulo Rob.-CAM-2NA, Guilherme-Eletropneumática-2NA, Guilherme-Eletro-Hidráulica-2NA]"
$ws.Range("D18").Value = "Andre B.-Circuitos Elétricos 2-"
$ws.Range("E18").Value = "[Paulo Rob.-CAM-2NA, Leandro-Sistemas de Refrigeração-2NA, -, -]"
$ws.Range("F18").Value = "[-, Allan Cupertino-Instalções Elétricas-2NA, -, -]"

# Row 19
$ws.Range("B19").Value = "[Allan Cupertino-Instalções Elétricas-2NA, Emerson-Eletrônica Básica-2NA]"
$ws.Range("C19").Value = "[Cláudio-Tecnologia da Soldagem-2NA, Paulo Rob.-CAM-2NA, Guilherme-Eletropneumática-2NA, Guilherme-Eletro-Hidráulica-2NA]"
$ws.Range("D19").Value = "Andre B.-Circuitos Elétricos 2-"
$ws.Range("E19").Value = "[Paulo Rob.-CAM-2NA, Leandro-Sistemas de Refrigeração-2NA, -, -]"
$ws.Range("F19").Value = "[-, Cleidson-Automação Industrial-2NA, -, -]"

# Row 20
$ws.Range("B20").Value = "[João Paulo-Lab. Circuitos Elétricos-2NA, Allan Cupertino-Instalções Elétricas-2NA]"
$ws.Range("C20").Value = "[Cláudio-Tecnologia da Soldagem-2NA, Leandro-Sistemas de Refrigeração-2NA, Guilherme-Eletropneumática-2NA, Guilherme-Eletro-Hidráulica-2NA]"
$ws.Range("D20").Value = "Allan Cupertino-Máquinas Elétricas-"
$ws.Range("E20").Value = "[Weslei-CAD-2NA, Weslei-CAD-2NA]"
$ws.Range("F20").Value = "[-, Cleidson-Automação Industrial-2NA, -, -]"

# Row 21
$ws.Range("B21").Value = "[Allan Cupertino-Lab. De Máquinas elétricas-2NA, Allan Cupertino-Lab. De Máquinas elétricas-2NA]"
$ws.Range("C21").Value = "[Cláudio-Tecnologia da Soldagem-2NA, Leandro-Sistemas de Refrigeração-2NA, Guilherme-Eletropneumática-2NA, Guilherme-Eletro-Hidráulica-2NA]"
$ws.Range("D21").Value = "Allan Cupertino-Máquinas Elétricas-"
$ws.Range("E21").Value = "[Weslei-CAD-2NA, Weslei-CAD-2NA]"
$ws.Range("F21").Value = "[-, Cleidson-Automação Industrial-2NA, -, -]"
